# Merge with 4-months-old stale changes on CDT laptop
#
# - Converts the "range" column (A2:A5) on the "A220-300" and "A320 neo"
#   sheets from hard-coded nautical-mile values into formulas that convert
#   them to kilometres (value * 1.852), keeping the original literal as the
#   formula operand.
# - Updates the active sheet / selected cell on each sheet: "ATR 72-600"
#   becomes the active tab with H13 selected, while "A220-300" and
#   "A320 neo" both end up with A6 selected (and are no longer the active
#   tab).

$wb = $excel.ActiveWorkbook

$wsATR  = $wb.Worksheets.Item("ATR 72-600")
$wsA220 = $wb.Worksheets.Item("A220-300")
$wsA320 = $wb.Worksheets.Item("A320 neo")

# --- A220-300: rewrite A2:A5 as formulas converting nm -> km -----------
$wsA220.Range("A2").Formula = "=0.159337248422758*1.852"
$wsA220.Range("A3").Formula = "=2146.44492328111*1.852"
$wsA220.Range("A4").Formula = "=3213.17696200406*1.852"
$wsA220.Range("A5").Formula = "=3860.43056465351*1.852"

# --- A320 neo: rewrite A2:A5 as formulas converting nm -> km -----------
$wsA320.Range("A2").Formula = "=-2.84575981787137*1.852"
$wsA320.Range("A3").Formula = "=2492.88560045532*1.852"
$wsA320.Range("A4").Formula = "=3497.43881616391*1.852"
$wsA320.Range("A5").Formula = "=4393.85315879339*1.852"

# --- Update selections / active sheet -----------------------------------
# Select A6 on the two data sheets first (leaving them not the active tab).
$wsA220.Activate()
$wsA220.Range("A6").Select()

$wsA320.Activate()
$wsA320.Range("A6").Select()

# Finally activate "ATR 72-600" and select H13, making it the active tab.
$wsATR.Activate()
$wsATR.Range("H13").Select()
